# Update "想去人数" (F column) counts across sheets, per commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 52
$ws1.Range("F4").Value = 1295
$ws1.Range("F6").Value = 372
$ws1.Range("F7").Value = 1194
$ws1.Range("F9").Value = 7352
$ws1.Range("F13").Value = 8016
$ws1.Range("F15").Value = 55
$ws1.Range("F16").Value = 5521
$ws1.Range("F17").Value = 52
$ws1.Range("F18").Value = 2438
$ws1.Range("F24").Value = 11
$ws1.Range("F25").Value = 388
$ws1.Range("F26").Value = 260
$ws1.Range("F28").Value = 2407
$ws1.Range("F30").Value = 270
$ws1.Range("F31").Value = 87
$ws1.Range("F33").Value = 590
$ws1.Range("F36").Value = 1519
$ws1.Range("F39").Value = 2364

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 98
$ws2.Range("F4").Value = 71
$ws2.Range("F5").Value = 12
$ws2.Range("F6").Value = 29

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 52
$ws4.Range("F6").Value = 1295
$ws4.Range("F7").Value = 98
$ws4.Range("F8").Value = 372
$ws4.Range("F9").Value = 1194
$ws4.Range("F11").Value = 7352
$ws4.Range("F15").Value = 8016
$ws4.Range("F17").Value = 55
$ws4.Range("F18").Value = 5521
$ws4.Range("F19").Value = 52
$ws4.Range("F20").Value = 2438
$ws4.Range("F26").Value = 11
$ws4.Range("F27").Value = 71
$ws4.Range("F28").Value = 388
$ws4.Range("F30").Value = 2407
$ws4.Range("F32").Value = 270
$ws4.Range("F33").Value = 87
$ws4.Range("F35").Value = 12
$ws4.Range("F36").Value = 590
$ws4.Range("F39").Value = 29
$ws4.Range("F40").Value = 1519
$ws4.Range("F43").Value = 2364
